$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A29").Value = "Record"
$ws.Range("B29").Value = "Balanço Geral"
$ws.Range("C29").Value = "Saúde"
$ws.Range("D29").Value = "2025-04-01T13:22"
$ws.Range("E29").Value = "Neutro"
$ws.Range("F29").Value = "Idosa mordida por cachorro teve que tomar vacina contra a raiva. *nota coberta*"
